$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.652.62'
$ws.Range("E2").Value = '  -5.77%  '
$ws.Range("D3").Value = '1.810.55'
$ws.Range("E3").Value = '  -4.96%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '275.76'
$ws.Range("E5").Value = '  -9.96%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("E7").Value = '  -6.62%  '
$ws.Range("D8").Value = '0.3493'
$ws.Range("E8").Value = '  -8.25%  '
$ws.Range("D9").Value = '44.16'
$ws.Range("E9").Value = '  -4.01%  '
$ws.Range("D10").Value = '0.06609'
$ws.Range("E10").Value = '  -9.34%  '
$ws.Range("D11").Value = '20.09'
$ws.Range("E11").Value = '  -9.68%  '
$ws.Range("D12").Value = '0.8385'
$ws.Range("E12").Value = '  -7.37%  '
$ws.Range("D13").Value = '0.07820'
$ws.Range("E13").Value = '  -4.74%  '
$ws.Range("D14").Value = '1.815.33'
$ws.Range("E14").Value = '  +68.02%  '
$ws.Range("D15").Value = '5.044'
$ws.Range("E15").Value = '  -5.59%  '
$ws.Range("D16").Value = '87.20'
$ws.Range("E16").Value = '  -8.98%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '13.87'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '0.000007945'
$ws.Range("E20").Value = '  -8.28%  '
$ws.Range("D21").Value = '25.756.22'
$ws.Range("E21").Value = '  -5.46%  '
$ws.Range("D22").Value = '4.709'
$ws.Range("E22").Value = '  -6.70%  '
$ws.Range("D23").Value = '9.992'
$ws.Range("E23").Value = '  -7.27%  '
$ws.Range("D24").Value = '6.066'
$ws.Range("E24").Value = '  -6.95%  '
$ws.Range("D25").Value = '141.21'
$ws.Range("E25").Value = '  -5.25%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '1.659'
$ws.Range("E26").Value = '  -5.03%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.085'
$ws.Range("E27").Value = '  -8.95%  '
$ws.Range("E28").Value = '  -8.51%  '
$ws.Range("E29").Value = '  -7.20%  '
$ws.Range("D30").Value = '4.306'
$ws.Range("E30").Value = '  -10.61%  '
$ws.Range("D31").Value = '4.198'
$ws.Range("E31").Value = '  -11.17%  '
$ws.Range("D32").Value = '0.08777'
$ws.Range("E32").Value = '  -4.84%  '
$ws.Range("D33").Value = '0.04848'
$ws.Range("E33").Value = '  -4.61%  '
$ws.Range("D34").Value = '0.7377'
$ws.Range("E34").Value = '  -11.01%  '
$ws.Range("D35").Value = '1.128'
$ws.Range("E35").Value = '  -7.38%  '
$ws.Range("D36").Value = '2.871'
$ws.Range("E36").Value = '  -4.49%  '
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = '3.039'
$ws.Range("E38").Value = '  -8.84%  '
$ws.Range("D39").Value = '2.430'
$ws.Range("E39").Value = '  -9.11%  '
$ws.Range("D40").Value = '0.5293'
$ws.Range("E40").Value = '  -9.36%  '
$ws.Range("E41").Value = '  -6.67%  '
$ws.Range("D42").Value = '0.9711'
$ws.Range("E42").Value = '  -9.80%  '
$ws.Range("D43").Value = '6.195'
$ws.Range("E43").Value = '  -6.51%  '
$ws.Range("D44").Value = '111.34'
$ws.Range("E44").Value = '  -4.99%  '
$ws.Range("D45").Value = '8.141'
$ws.Range("E45").Value = '  -12.58%  '
$ws.Range("D46").Value = '0.4654'
$ws.Range("E46").Value = '  -8.47%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").Value = '0.1384'
$ws.Range("E48").Value = '  -9.10%  '
$ws.Range("D49").Value = '9.161'
$ws.Range("E49").Value = '  -9.38%  '
$ws.Range("D50").Value = '35.57'
$ws.Range("E50").Value = '  -7.37%  '
$ws.Range("D51").Value = '1.489'
$ws.Range("E51").Value = '  -9.24%  '
